$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("startup")

# Update the "FilesTab" query text in B4: drop the `File Type` return column
# and the `Breed` return column from the Cypher query.
$newQuery = "MATCH (f:file)-->(parent)`nWITH DISTINCT f, parent`nMATCH (f)-[*]->(c:case)<--(demo:demographic)`n MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)`n MATCH (samp:sample)-->(c) `n WHERE samp.specific_sample_pathology IN [`"Osteosarcoma`"]  `nWITH DISTINCT f, parent, c, demo, diag, s`nRETURN coalesce(f.file_name, '') AS ``File Name``, `n        coalesce(labels(parent)[0], '') AS ``Association``,`n        coalesce(f.file_description, '') AS ``Description``,`n        coalesce(f.file_format, '') AS ``Format``,`n        coalesce(f.file_size, '') AS ``Size``,`n        coalesce(c.case_id, '') AS ``Case ID``, `n        coalesce(diag.disease_term,'') AS Diagnosis , `n        coalesce(s.clinical_study_designation,'') AS ``Study Code``"

$ws.Range("B4").Value = $newQuery

# Move the saved selection/active cell on the sheet from D8 to B4.
$ws.Range("B4").Select()
